$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4 (shifting the existing "fin" row down to row 5)
$ws.Rows.Item(4).Insert()

# Fill the new row 4 with data matching pattern of row 3 (python_bases / PYB-100 / Types composés)
$ws.Range("A4").Value = "python_bases"
$ws.Range("B4").Value = "PYB-100"
$ws.Range("C4").Value = "Types composés"

# Update selection to match diff (C5 selected)
$ws.Range("C5").Select()
